$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 794.875
$ws.Range("I103").Value = 474
$ws.Range("J103").Value = 987.4
$ws.Range("K103").Value = 1422
$ws.Range("L103").Value = 2962.2
$ws.Range("M103").Value = -836
$ws.Range("N103").Value = -4134.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 51454.95
$ws.Range("I2").Value = 63838.062
$ws.Range("J2").Value = 1922.5
$ws.Range("K2").Value = 63838.062
$ws.Range("L2").Value = 1922.5
$ws.Range("M2").Value = -63725.062
$ws.Range("N2").Value = -2148.5
$ws.Range("H32").Value = 20723.246
$ws.Range("I32").Value = 3324.22
$ws.Range("J32").Value = 145002
$ws.Range("K32").Value = 3324.22
$ws.Range("L32").Value = 145002
$ws.Range("M32").Value = -3037.22
$ws.Range("N32").Value = -145576
$ws.Range("H61").Value = 2196.8438
$ws.Range("I61").Value = 1657.16
$ws.Range("J61").Value = 4124.2856
$ws.Range("K61").Value = 1657.16
$ws.Range("L61").Value = 4124.2856
$ws.Range("M61").Value = -1445.16
$ws.Range("N61").Value = -4548.2856
$ws.Range("H74").Value = 7587.3823
$ws.Range("I74").Value = 1192.6129
$ws.Range("K74").Value = 1192.6129
$ws.Range("M74").Value = -318.6129000000001
$ws.Range("H77").Value = 7587.3823
$ws.Range("I77").Value = 1192.6129
$ws.Range("K77").Value = 5963.0645
$ws.Range("M77").Value = -1595.0645
$ws.Range("H102").Value = 3598
$ws.Range("I102").Value = 5995
$ws.Range("K102").Value = 5995
$ws.Range("M102").Value = -4373
$ws.Range("H116").Value = 51454.95
$ws.Range("I116").Value = 63838.062
$ws.Range("J116").Value = 1922.5
$ws.Range("K116").Value = 63838.062
$ws.Range("L116").Value = 1922.5
$ws.Range("M116").Value = -61544.062
$ws.Range("N116").Value = -6510.5
$ws.Range("H132").Value = 4552
$ws.Range("I132").Value = 4516.1055
$ws.Range("J132").Value = 4722.5
$ws.Range("K132").Value = 13548.3165
$ws.Range("L132").Value = 14167.5
$ws.Range("M132").Value = -11018.3165
$ws.Range("N132").Value = -19227.5
$ws.Range("H136").Value = 2196.8438
$ws.Range("I136").Value = 1657.16
$ws.Range("J136").Value = 4124.2856
$ws.Range("K136").Value = 4971.48
$ws.Range("L136").Value = 12372.8568
$ws.Range("M136").Value = -2421.48
$ws.Range("N136").Value = -17472.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 51454.95
$ws.Range("I3").Value = 63838.062
$ws.Range("J3").Value = 1922.5
$ws.Range("K3").Value = 63838.062
$ws.Range("L3").Value = 1922.5
$ws.Range("M3").Value = -63724.062
$ws.Range("N3").Value = -2150.5
$ws.Range("H86").Value = 8246.25
$ws.Range("I86").Value = 2146.4546
$ws.Range("J86").Value = 21665.8
$ws.Range("K86").Value = 2146.4546
$ws.Range("L86").Value = 21665.8
$ws.Range("M86").Value = -1023.4546
$ws.Range("N86").Value = -23911.8
$ws.Range("H89").Value = 8246.25
$ws.Range("I89").Value = 2146.4546
$ws.Range("J89").Value = 21665.8
$ws.Range("K89").Value = 10732.273
$ws.Range("L89").Value = 108329
$ws.Range("M89").Value = -5116.273000000001
$ws.Range("N89").Value = -119561
$ws.Range("H94").Value = 1162.762
$ws.Range("I94").Value = 1235.5
$ws.Range("J94").Value = 1017.2857
$ws.Range("K94").Value = 1235.5
$ws.Range("L94").Value = 1017.2857
$ws.Range("M94").Value = -784.5
$ws.Range("N94").Value = -1919.2857
$ws.Range("H134").Value = 2903.1707
$ws.Range("I134").Value = 1810.625
$ws.Range("J134").Value = 4445.5884
$ws.Range("K134").Value = 5431.875
$ws.Range("L134").Value = 13336.7652
$ws.Range("M134").Value = -2896.875
$ws.Range("N134").Value = -18406.7652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1747.7916
$ws.Range("I122").Value = 1094.7693
$ws.Range("J122").Value = 2519.5454
$ws.Range("K122").Value = 3284.3079
$ws.Range("L122").Value = 7558.6362
$ws.Range("M122").Value = -834.3078999999998
$ws.Range("N122").Value = -12458.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2148.389
$ws.Range("I5").Value = 2279.4285
$ws.Range("J5").Value = 2065
$ws.Range("K5").Value = 6838.2855
$ws.Range("L5").Value = 6195
$ws.Range("M5").Value = -6726.2855
$ws.Range("N5").Value = -6419
$ws.Range("H131").Value = 1479.9138
$ws.Range("I131").Value = 373
$ws.Range("J131").Value = 1683.2245
$ws.Range("K131").Value = 1119
$ws.Range("L131").Value = 5049.6735
$ws.Range("M131").Value = 3921
$ws.Range("N131").Value = -15129.6735
$ws.Range("H135").Value = 2148.389
$ws.Range("I135").Value = 2279.4285
$ws.Range("J135").Value = 2065
$ws.Range("K135").Value = 20514.8565
$ws.Range("L135").Value = 18585
$ws.Range("M135").Value = -17979.8565
$ws.Range("N135").Value = -23655

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12380
$ws.Range("I70").Value = 20450
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 20450
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -20180
$ws.Range("N70").Value = -7540
$ws.Range("H73").Value = 12380
$ws.Range("I73").Value = 20450
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 20450
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = -19514
$ws.Range("N73").Value = -8872
$ws.Range("H97").Value = 1702
$ws.Range("I97").Value = 1655
$ws.Range("J97").Value = 1733.3334
$ws.Range("K97").Value = 1655
$ws.Range("L97").Value = 1733.3334
$ws.Range("M97").Value = -1159
$ws.Range("N97").Value = -2725.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2776.4062
$ws.Range("I7").Value = 1812.9166
$ws.Range("K7").Value = 1812.9166
$ws.Range("M7").Value = -1700.9166
$ws.Range("H126").Value = 2776.4062
$ws.Range("I126").Value = 1812.9166
$ws.Range("K126").Value = 5438.7498
$ws.Range("M126").Value = -2968.7498
$ws.Range("H136").Value = 5221.5757
$ws.Range("I136").Value = 3199.875
$ws.Range("J136").Value = 10612.777
$ws.Range("K136").Value = 9599.625
$ws.Range("L136").Value = 31838.331
$ws.Range("M136").Value = -7049.625
$ws.Range("N136").Value = -36938.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8198185.5
$ws.Range("I132").Value = 11906036
$ws.Range("J132").Value = 1885.5264
$ws.Range("K132").Value = 35718108
$ws.Range("L132").Value = 5656.5792
$ws.Range("M132").Value = -35715578
$ws.Range("N132").Value = -10716.5792
